$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new row 19 (sc14, the "small" version) first, so the new
# shared-string entries are interned in the same order as the target file:
# "...- small" (E19) then "...- big" (E18) then "sc14" (A19).
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = "A(5803) Single hierarchy conditions, two policies. one extra condition on one policy - small"
$ws.Range("F19").Value = "input files"
$ws.Range("G19").Value = "no"
$ws.Range("H19").Value = "to do"

# Update row 18 (sc13) - now becomes the "big" version
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = "A(5803) Single hierarchy conditions, two policies. one extra condition on one policy - big"

# Set A19 last so "sc14" is interned after the small/big description strings
$ws.Range("A19").Value = "sc14"

# Adjust column E width to fit the new, longer description text
# (closest achievable value to the target 75.6640625 char-width units)
$ws.Columns("E").ColumnWidth = 74.8

# Update selection to C19
$ws.Range("C19").Select()
